$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: weekday numbers -> weekday names
$ws.Range("B1").Value = "segunda"
$ws.Range("C1").Value = "terça"
$ws.Range("D1").Value = "quarta"
$ws.Range("E1").Value = "quinta"
$ws.Range("F1").Value = "sexta"

# Column A: period index -> start time of each class slot
$ws.Range("A2").Value = "7:00"
$ws.Range("A3").Value = "7:50"
$ws.Range("A4").Value = "8:40"
$ws.Range("A5").Value = "9:30"
$ws.Range("A6").Value = "10:40"
$ws.Range("A7").Value = "11:30"
$ws.Range("A8").Value = "13:00"
$ws.Range("A9").Value = "13:50"
$ws.Range("A10").Value = "14:40"
$ws.Range("A11").Value = "15:30"
$ws.Range("A12").Value = "16:40"
$ws.Range("A13").Value = "17:30"

# Move "Desenho Técnico" from D3 to B10
$ws.Range("D3").Value = "-"
$ws.Range("B10").Value = "Desenho Técnico"

# Move "EAP" from F3 to F2
$ws.Range("F3").Value = "-"
$ws.Range("F2").Value = "EAP"

# Move "EAP" from E13 to F11
$ws.Range("E13").Value = "-"
$ws.Range("F11").Value = "EAP"

# Move "Circuitos Elétricos 2" from B12 to B13
$ws.Range("B12").Value = "-"
$ws.Range("B13").Value = "Circuitos Elétricos 2"
